# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.813.46"
$ws.Range("E2").Value = "  +6.74%  "
$ws.Range("D3").Value = "2.427.41"
$ws.Range("E3").Value = "  +6.95%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'477.80"
$ws.Range("E5").Value = "  +10.62%  "
$ws.Range("D6").Value = "'138.96"
$ws.Range("E6").Value = "  +20.85%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.503"
$ws.Range("E8").Value = "  +11.85%  "
$ws.Range("D9").Value = "2.451.22"
$ws.Range("E9").Value = "  +8.21%  "
$ws.Range("D10").Value = "'0.0958"
$ws.Range("E10").Value = "  +14.60%  "
$ws.Range("D11").Value = "'5.46"
$ws.Range("E11").Value = "  +7.67%  "
$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = "  +10.34%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "2.855.92"
$ws.Range("E14").Value = "  +8.50%  "
$ws.Range("D15").Value = "54.962.30"
$ws.Range("E15").Value = "  +6.91%  "
$ws.Range("D16").Value = "'20.41"
$ws.Range("E16").Value = "  +12.71%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("E17").Value = "  +19.58%  "
$ws.Range("D18").Value = "2.449.77"
$ws.Range("E18").Value = "  +8.90%  "
$ws.Range("D19").Value = "'4.33"
$ws.Range("E19").Value = "  +12.75%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'313.01"
$ws.Range("E20").Value = "  +8.72%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'9.84"
$ws.Range("E21").Value = "  +17.66%  "
$ws.Range("D22").Value = "'0.993"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'5.62"
$ws.Range("E23").Value = "  +14.72%  "
$ws.Range("D24").Value = "'57.07"
$ws.Range("E24").Value = "  +8.93%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'0.400"
$ws.Range("E26").Value = "  +11.96%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  +17.99%  "
$ws.Range("D28").Value = "2.541.44"
$ws.Range("E28").Value = "  +9.53%  "
$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  +12.10%  "
$ws.Range("D30").Value = "0.0₃0765"
$ws.Range("E30").Value = "  +24.79%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'148.54"
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("D33").Value = "'17.85"
$ws.Range("E33").Value = "  +9.99%  "
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = "  +13.53%  "
$ws.Range("D35").Value = "'5.14"
$ws.Range("E35").Value = "  +13.66%  "
$ws.Range("E36").Value = "  +18.15%  "
$ws.Range("D37").Value = "'3.57"
$ws.Range("E37").Value = "  +9.62%  "
$ws.Range("D38").Value = "'0.836"
$ws.Range("E38").Value = "  +15.27%  "
$ws.Range("D39").Value = "'33.47"
$ws.Range("E39").Value = "  +6.32%  "
$ws.Range("D40").Value = "'0.992"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'0.600"
$ws.Range("E41").Value = "  +9.57%  "
$ws.Range("D42").Value = "'3.40"
$ws.Range("E42").Value = "  +12.97%  "
$ws.Range("D43").Value = "'0.0541"
$ws.Range("E43").Value = "  +12.56%  "
$ws.Range("D44").Value = "'1.28"
$ws.Range("E44").Value = "  +17.09%  "
$ws.Range("D45").Value = "'10.14"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'4.61"
$ws.Range("E46").Value = "  +21.83%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.0891"
$ws.Range("E47").Value = "  +14.42%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'250.90"
$ws.Range("E48").Value = "  +32.95%  "
$ws.Range("D49").Value = "'0.0221"
$ws.Range("E49").Value = "  +12.77%  "
$ws.Range("D50").Value = "1.917.34"
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("D51").Value = "'17.14"
$ws.Range("E51").Value = "  +13.92%  "
